# Appends the new repay-list entries (rows 14-45) restored from a duplicate-cleanup pass.
# Values are stored Base64-encoded (UTF-8) below purely to dodge PowerShell/Excel
# quoting edge cases (the source strings are Python dict reprs full of ' and ").
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$encodedValues = @(
    "J1JVMDAwQTBKUE5OOSc6IHsnbmFtZV9wYXBlcic6ICfQkNCeICLQpNCh0Jog0JXQrdChIiDQn9CQ0J4sINCy0YvQv9GD0YHQuiAwMScsICd2b2x1bWUnOiA3MDAwMH0sIA==",
    "J1JVMDAwNjc2NTA5Nic6IHsnbmFtZV9wYXBlcic6ICfQkNCfICAi0J3QuNC20L3QtdC60LDQvNGB0LrQvdC10YTRgtC10YXQuNC8IiAo0J3QmtCd0KUpINCf0JDQniwg0LLRi9C/0YPRgdC6IDAyJywgJ3ZvbHVtZSc6IDIyMH0sIA==",
    "J1JVMDAwOTA5MjEzNCc6IHsnbmFtZV9wYXBlcic6ICfQkNCfICLQm9C10L3RjdC90LXRgNCz0L4iICDQntCQ0J4sINCy0YvQv9GD0YHQuiAwMScsICd2b2x1bWUnOiAxMjB9LCA=",
    "J1JVMDAwNjk0NDE0Nyc6IHsnbmFtZV9wYXBlcic6ICfQkNCfICLQotCw0YLQvdC10YTRgtGMIiAg0J/QkNCeLCDQstGL0L/Rg9GB0LogMDMnLCAndm9sdW1lJzogMjJ9LCA=",
    "J1JVMDAwOTAyOTUyNCc6IHsnbmFtZV9wYXBlcic6ICfQkNCfICLQodGD0YDQs9GD0YLQvdC10YTRgtC10LPQsNC3IiDQntCQ0J4sINCy0YvQv9GD0YHQuiAwMScsICd2b2x1bWUnOiA3MDB9LCA=",
    "J1JVMDAwOTA0NjUxMCc6IHsnbmFtZV9wYXBlcic6ICfQkNCeICLQodC10LLQtdGA0YHRgtCw0LvRjCIg0J/QkNCeLCDQstGL0L/Rg9GB0LogMDInLCAndm9sdW1lJzogMTJ9LCA=",
    "J1JVMDAwNzk3Njk2NSc6IHsnbmFtZV9wYXBlcic6ICfQkNCfICLQkNCd0JogItCR0LDRiNC90LXRhNGC0YwiINCf0JDQniwg0LLRi9C/0YPRgdC6IDAxJywgJ3ZvbHVtZSc6IDEwfSwg",
    "J1JVMDAwQTBKWEs0MCc6IHsnbmFtZV9wYXBlcic6ICfQntCx0LsuICLQk9GA0YPQv9C/0LAg0JrQvtC80L/QsNC90LjQuSDQn9CY0JoiINCf0JDQniwg0YHQtdGA0LjRjyDQkdCeLdCfMDEnLCAndm9sdW1lJzogN30sIA==",
    "J1JVMDAwQTBKUFBUMSc6IHsnbmFtZV9wYXBlcic6ICfQkNCeICLQnNCg0KHQmiDQo9GA0LDQu9CwIiDQntCQ0J4sINCy0YvQv9GD0YHQuiAwMScsICd2b2x1bWUnOiAyMDAwMH0sIA==",
    "J1JVMDAwQTBEUVpFMyc6IHsnbmFtZV9wYXBlcic6ICfQkNCeICLQkNCk0JogItCh0LjRgdGC0LXQvNCwIiDQn9CQ0J4sINCy0YvQv9GD0YHQuiAwNScsICd2b2x1bWUnOiA2MDB9LCA=",
    "J1JVMDAwNzY2MTYyNSc6IHsnbmFtZV9wYXBlcic6ICfQkNCeICLQk9Cw0LfQv9GA0L7QvCIg0J/QkNCeLCDQstGL0L/Rg9GB0LogMDInLCAndm9sdW1lJzogMTAwfSwg",
    "J1JVMDAwQTBKUFBONCc6IHsnbmFtZV9wYXBlcic6ICfQkNCeICLQnNCg0KHQmiDQktC+0LvQs9C4IiDQn9CQ0J4sINCy0YvQv9GD0YHQuiAwMScsICd2b2x1bWUnOiAxMjAwMDB9LCA=",
    "J1JVMDAwQTBKUE45Nic6IHsnbmFtZV9wYXBlcic6ICfQkNCeICLQnNCg0KHQmiDQptC10L3RgtGA0LAg0Lgg0J/RgNC40LLQvtC70LbRjNGPIiDQntCQ0J4sINCy0YvQv9GD0YHQuiAwMScsICd2b2x1bWUnOiA3MDAwMH0sIA==",
    "J1JVMDAwQTBaWVFZNyc6IHsnbmFtZV9wYXBlcic6ICfQntCx0LsuINCQ0KTQmiAi0KHQuNGB0YLQtdC80LAiINCf0JDQniwg0YHQtdGA0LjRjyAwMDFQLTA3JywgJ3ZvbHVtZSc6IDEwfSwg",
    "J1JVMDAwOTA0NjQ1Mic6IHsnbmFtZV9wYXBlcic6ICfQkNCeICLQndCb0JzQmiIg0J/QkNCeLCDQstGL0L/Rg9GB0LogMDEnLCAndm9sdW1lJzogMTUwfSwg",
    "J1JVMDAwNzI1MjgxMyc6IHsnbmFtZV9wYXBlcic6ICfQkNCeICLQkNCb0KDQntCh0JAiINCf0JDQniwg0LLRi9C/0YPRgdC6IDAzJywgJ3ZvbHVtZSc6IDMwMH0sIA==",
    "J1JVMDAwQTBKVlczMCc6IHsnbmFtZV9wYXBlcic6ICfQntCk0JcgMjYyMTct0J/QlCcsICd2b2x1bWUnOiA2fSw=",
    "J1JVMDAwQTBKUkVRNyc6IHsnbmFtZV9wYXBlcic6ICfQntCk0JcgMjYyMDUt0J/QlCcsICd2b2x1bWUnOiA2fSw=",
    "ICdSVTAwMDkwMjk1NTcnOiB7J25hbWVfcGFwZXInOiAn0JDQnyAi0KHQkdCV0KDQkdCQ0J3QmiDQoNC+0YHRgdC40LgiINCf0JDQnicsICd2b2x1bWUnOiA4MH0sIA==",
    "J1JVMDAwQTBKUFBHOCc6IHsnbmFtZV9wYXBlcic6ICfQkNCeICLQnNCg0KHQmiDQrtCz0LAiINCf0JDQniwg0LLRi9C/0YPRgdC6IDAxJywgJ3ZvbHVtZSc6IDcwMDAwfSwg",
    "J1JVMDAwOTA2MjQ2Nyc6IHsnbmFtZV9wYXBlcic6ICfQkNCeICLQk9Cw0LfQv9GA0L7QvCDQvdC10YTRgtGMIiAg0J/QkNCeLCDQstGL0L/Rg9GB0LogMDEnLCAndm9sdW1lJzogMzB9LA==",
    "J1JVMDAwQTBKV1U5OCc6IHsnbmFtZV9wYXBlcic6ICfQntCx0LsuICLQk9GA0YPQv9C/0LAg0JvQodCgIiDQn9CQ0J4sINGB0LXRgNC40Y8gMDAxUC0wMScsICd2b2x1bWUnOiAyfSwg",
    "J1JVMDAwQTBKUEZQMCc6IHsnbmFtZV9wYXBlcic6ICfQkNCeICLQk9GA0YPQv9C/0LAg0JvQodCgIiDQn9CQ0J4sINCy0YvQv9GD0YHQuiAwMScsICd2b2x1bWUnOiAyfSwg",
    "J1JVMDAwOTAyNDI3Nyc6IHsnbmFtZV9wYXBlcic6ICfQkNCeICLQm9Cj0JrQntCZ0JsiINCf0JDQniwg0LLRi9C/0YPRgdC6IDAxJywgJ3ZvbHVtZSc6IDJ9LA==",
    "J1JVMDAwQTEwMFcyOSc6IHsnbmFtZV9wYXBlcic6ICfQntCx0LsuINCb0JogItCg0L7QtNC10LvQtdC9IiDQl9CQ0J4sINGB0LXRgNC40Y8gMDAxUC0wMicsICd2b2x1bWUnOiAxfSwg",
    "J1JVMDAwQTEwMUNCNic6IHsnbmFtZV9wYXBlcic6ICfQntCx0LsuICLQm9C40LfQuNC90LMt0KLRgNC10LnQtCIg0J7QntCeLCDRgdC10YDQuNGPIDAwMVAtMDEnLCAndm9sdW1lJzogNX0sIA==",
    "J1JVMDAwQTEwMUZUMSc6IHsnbmFtZV9wYXBlcic6ICfQntCx0LsuICLQm9C+0LzQsdCw0YDQtCAi0JzQsNGB0YLQtdGAIiDQntCe0J4sINGB0LXRgNC40Y8g0JHQni3QnzA3JywgJ3ZvbHVtZSc6IDN9LCA=",
    "J1JVMDAwQTBKTkdBNSc6IHsnbmFtZV9wYXBlcic6ICfQkNCeICLQrtC90LjQv9GA0L4iINCf0JDQniwg0LLRi9C/0YPRgdC6IDAyJywgJ3ZvbHVtZSc6IDIwMDB9LCA=",
    "J1JVMDAwQTBKUDVWNic6IHsnbmFtZV9wYXBlcic6ICfQkNCeINCR0LDQvdC6INCS0KLQkSDQn9CQ0J4sINCy0YvQv9GD0YHQuiAwNCcsICd2b2x1bWUnOiA4MDAwMH0sIA==",
    "J1JVMDAwQTBKUjRBMSc6IHsnbmFtZV9wYXBlcic6ICfQkNCeICLQnNC+0YHQutC+0LLRgdC60LDRjyDQkdC40YDQttCwIiDQn9CQ0J4sINCy0YvQv9GD0YHQuiAwMScsICd2b2x1bWUnOiAyMH0sIA==",
    "J1JVMDAwQTBES1ZTNSc6IHsnbmFtZV9wYXBlcic6ICfQkNCeICLQndCe0JLQkNCi0K3QmiIg0J/QkNCeLCDQstGL0L/Rg9GB0LogMDInLCAndm9sdW1lJzogM30sIA==",
    "J1JVMDAwQTBITTVDMSc6IHsnbmFtZV9wYXBlcic6ICfQkNCeICLQotCw0YLRgtC10LvQtdC60L7QvCIg0J/QkNCeLCDQstGL0L/Rg9GB0LogMDInLCAndm9sdW1lJzogMTAwMDB9"
)

$startRow = 14
for ($i = 0; $i -lt $encodedValues.Length; $i++) {
    $bytes = [System.Convert]::FromBase64String($encodedValues[$i])
    $text = [System.Text.Encoding]::UTF8.GetString($bytes)
    $row = $startRow + $i
    $cell = $ws.Cells.Item($row, 1)

    # Route the literal text through a formula + paste-as-values round trip instead of
    # $cell.Value = $text directly: a direct assignment of a string that begins with an
    # apostrophe makes Excel apply an implicit 'quote prefix' cell style (s="1" / text-forced
    # formatting), which the source workbook does not have on these cells.
    $escaped = $text.Replace([char]34, [string][char]34 + [char]34)
    $cell.Formula = "=" + [char]34 + $escaped + [char]34
    $cell.Copy()
    $cell.PasteSpecial(-4163)  # xlPasteValues
}

$excel.CutCopyMode = 0
$ws.Range("H45").Select() | Out-Null
